$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "CreatedUser" column (J). This shifts the "Accessory"
# column (K) left into J, matching the bug fix for the Feature Master import.
[void]$ws.Columns("J").Select()
$ws.Columns("J").Delete()
